$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 66 (shifts the
# existing rows 66-109 down to 68-111, carrying along the date style ("s=2")
# already used by column D in every data row).
$ws.Range("66:67").Insert()

# Row 66: brand-new "Especial" price quote dated 44673.
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44673
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100108
$ws.Range("H66").Value = "Tropicales y subtropicales"
$ws.Range("I66").Value = 100108003
$ws.Range("J66").Value = "Maracuyá"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Especial"
$ws.Range("M66").Value = 120
$ws.Range("N66").Value = 24000
$ws.Range("O66").Value = 25000
$ws.Range("P66").Value = 24500
$ws.Range("Q66").Value = "$/caja 20 kilos"
$ws.Range("R66").Value = "Región de Arica y Parinacota"
$ws.Range("S66").Value = 1225
$ws.Range("T66").Value = 20

# Row 67: brand-new "Primera" price quote dated 44673.
$ws.Range("A67").Value = 1
$ws.Range("B67").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C67").Value = "Arica y Parinacota"
$ws.Range("D67").Value = 44673
$ws.Range("E67").Value = 15
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100108
$ws.Range("H67").Value = "Tropicales y subtropicales"
$ws.Range("I67").Value = 100108003
$ws.Range("J67").Value = "Maracuyá"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 130
$ws.Range("N67").Value = 21000
$ws.Range("O67").Value = 22000
$ws.Range("P67").Value = 21500
$ws.Range("Q67").Value = "$/caja 20 kilos"
$ws.Range("R67").Value = "Región de Arica y Parinacota"
$ws.Range("S67").Value = 1075
$ws.Range("T67").Value = 20
